{"js": "\n// Replace the 100 arithmetic-fact answers in the single 20x5 table, in\n// document order (row-major). The mapping is positional because some\n// \"before\" values repeat (e.g. \"87-69=18\" appears twice but maps to two\n// different \"after\" values), so a global text find/replace would be\n// ambiguous \u2014 we must address each table cell directly.\nconst OLD_VALUES = [\"25+5=30\", \"22+20=42\", \"99-73=26\", \"58+35=93\", \"64-29=35\", \"87-35=52\", \"58-10=48\", \"77-43=34\", \"58+6=64\", \"68-17=51\", \"83-77=6\", \"29-26=3\", \"5+1=6\", \"11+7=18\", \"1+51=52\", \"63-2=61\", \"0+11=11\", \"64+29=93\", \"83-6=77\", \"76-56=20\", \"2+60=62\", \"3+9=12\", \"2+20=22\", \"9+42=51\", \"33+8=41\", \"93-42=51\", \"82-3=79\", \"77-39=38\", \"90-30=60\", \"95+4=99\", \"29+46=75\", \"16+82=98\", \"26-24=2\", \"33-19=14\", \"10+31=41\", \"37+9=46\", \"12+33=45\", \"36-21=15\", \"51+1=52\", \"81-67=14\", \"60+22=82\", \"41-18=23\", \"73+24=97\", \"10+69=79\", \"94-46=48\", \"85-36=49\", \"74-43=31\", \"57-18=39\", \"86-4=82\", \"56-20=36\", \"38-19=19\", \"99-22=77\", \"95-30=65\", \"60-8=52\", \"44-36=8\", \"6+87=93\", \"59+40=99\", \"11+54=65\", \"59+6=65\", \"52-25=27\", \"19+57=76\", \"87-69=18\", \"85-58=27\", \"92-7=85\", \"41-17=24\", \"32+18=50\", \"79-11=68\", \"78-10=68\", \"99-70=29\", \"45-33=12\", \"49+10=59\", \"31-1=30\", \"33+29=62\", \"52-30=22\", \"51+32=83\", \"91-91=0\", \"58-14=44\", \"68+28=96\", \"39-37=2\", \"35+22=57\", \"72-4=68\", \"52+10=62\", \"7+68=75\", \"82-37=45\", \"5+49=54\", \"87-59=28\", \"82+17=99\", \"13+69=82\", \"14+85=99\", \"45-44=1\", \"26+41=67\", \"67+31=98\", \"14+23=37\", \"40+48=88\", \"36+58=94\", \"82-38=44\", \"27+33=60\", \"23-22=1\", \"87-69=18\", \"83-23=60\"];\nconst NEW_VALUES = [\"25+72=97\", \"77-30=47\", \"98-7=91\", \"72-51=21\", \"20+60=80\", \"53+46=99\", \"39-27=12\", \"26+23=49\", \"32+24=56\", \"28-9=19\", \"47+24=71\", \"15+64=79\", \"56-27=29\", \"36+38=74\", \"80-57=23\", \"60+2=62\", \"1+19=20\", \"36-36=0\", \"7+12=19\", \"28+1=29\", \"15-9=6\", \"81-6=75\", \"19+12=31\", \"74-59=15\", \"97-58=39\", \"44+54=98\", \"66+16=82\", \"1+87=88\", \"96-7=89\", \"24-13=11\", \"66-16=50\", \"8+56=64\", \"20+53=73\", \"10+29=39\", \"82-8=74\", \"46+45=91\", \"42+6=48\", \"98-88=10\", \"29+29=58\", \"28-12=16\", \"9+85=94\", \"60+16=76\", \"39-2=37\", \"89-19=70\", \"1+7=8\", \"10+63=73\", \"82-48=34\", \"86-7=79\", \"88-56=32\", \"8+63=71\", \"35+21=56\", \"4+25=29\", \"43-34=9\", \"18+77=95\", \"49+14=63\", \"79-50=29\", \"15+80=95\", \"8+54=62\", \"26+42=68\", \"34-28=6\", \"3+29=32\", \"65-29=36\", \"30+28=58\", \"69-14=55\", \"13+36=49\", \"88-56=32\", \"84-28=56\", \"61-24=37\", \"50+21=71\", \"34+21=55\", \"59-18=41\", \"86-52=34\", \"80-20=60\", \"3-2=1\", \"81-35=46\", \"45+38=83\", \"5+28=33\", \"87+12=99\", \"58+33=91\", \"37-4=33\", \"41+55=96\", \"32-32=0\", \"19+64=83\", \"43-2=41\", \"50+32=82\", \"4+22=26\", \"86+4=90\", \"88-53=35\", \"71-5=66\", \"34+54=88\", \"14+65=79\", \"38-2=36\", \"79-22=57\", \"40+47=87\", \"82-4=78\", \"67-32=35\", \"31-13=18\", \"27+31=58\", \"52-45=7\", \"26+56=82\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst COLS = 5;\nconst rowCount = table.rowCount;\nconst cells = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < COLS; c++) {\n    cells.push(table.getCell(r, c));\n  }\n}\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nif (cells.length !== OLD_VALUES.length) {\n  throw new Error(\n    \"Expected \" + OLD_VALUES.length + \" table cells, found \" + cells.length\n  );\n}\n\nfor (let i = 0; i < cells.length; i++) {\n  const cell = cells[i];\n  const expected = OLD_VALUES[i];\n  const actual = (cell.value || \"\").trim();\n  if (actual !== expected) {\n    throw new Error(\n      \"Cell \" + i + \": expected '\" + expected + \"' but found '\" + actual + \"'\"\n    );\n  }\n  cell.value = NEW_VALUES[i];\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-fact answers in the single 20x5 table, in\n# document order (row-major). The mapping is positional because some\n# \"before\" values repeat (e.g. \"87-69=18\" appears twice but maps to two\n# different \"after\" values), so a global Find/Replace would be ambiguous \u2014\n# each table cell is addressed directly by (row, col).\n\n$d = $word.ActiveDocument\n\n$oldValues = @(\"25+5=30\", \"22+20=42\", \"99-73=26\", \"58+35=93\", \"64-29=35\", \"87-35=52\", \"58-10=48\", \"77-43=34\", \"58+6=64\", \"68-17=51\", \"83-77=6\", \"29-26=3\", \"5+1=6\", \"11+7=18\", \"1+51=52\", \"63-2=61\", \"0+11=11\", \"64+29=93\", \"83-6=77\", \"76-56=20\", \"2+60=62\", \"3+9=12\", \"2+20=22\", \"9+42=51\", \"33+8=41\", \"93-42=51\", \"82-3=79\", \"77-39=38\", \"90-30=60\", \"95+4=99\", \"29+46=75\", \"16+82=98\", \"26-24=2\", \"33-19=14\", \"10+31=41\", \"37+9=46\", \"12+33=45\", \"36-21=15\", \"51+1=52\", \"81-67=14\", \"60+22=82\", \"41-18=23\", \"73+24=97\", \"10+69=79\", \"94-46=48\", \"85-36=49\", \"74-43=31\", \"57-18=39\", \"86-4=82\", \"56-20=36\", \"38-19=19\", \"99-22=77\", \"95-30=65\", \"60-8=52\", \"44-36=8\", \"6+87=93\", \"59+40=99\", \"11+54=65\", \"59+6=65\", \"52-25=27\", \"19+57=76\", \"87-69=18\", \"85-58=27\", \"92-7=85\", \"41-17=24\", \"32+18=50\", \"79-11=68\", \"78-10=68\", \"99-70=29\", \"45-33=12\", \"49+10=59\", \"31-1=30\", \"33+29=62\", \"52-30=22\", \"51+32=83\", \"91-91=0\", \"58-14=44\", \"68+28=96\", \"39-37=2\", \"35+22=57\", \"72-4=68\", \"52+10=62\", \"7+68=75\", \"82-37=45\", \"5+49=54\", \"87-59=28\", \"82+17=99\", \"13+69=82\", \"14+85=99\", \"45-44=1\", \"26+41=67\", \"67+31=98\", \"14+23=37\", \"40+48=88\", \"36+58=94\", \"82-38=44\", \"27+33=60\", \"23-22=1\", \"87-69=18\", \"83-23=60\")\n$newValues = @(\"25+72=97\", \"77-30=47\", \"98-7=91\", \"72-51=21\", \"20+60=80\", \"53+46=99\", \"39-27=12\", \"26+23=49\", \"32+24=56\", \"28-9=19\", \"47+24=71\", \"15+64=79\", \"56-27=29\", \"36+38=74\", \"80-57=23\", \"60+2=62\", \"1+19=20\", \"36-36=0\", \"7+12=19\", \"28+1=29\", \"15-9=6\", \"81-6=75\", \"19+12=31\", \"74-59=15\", \"97-58=39\", \"44+54=98\", \"66+16=82\", \"1+87=88\", \"96-7=89\", \"24-13=11\", \"66-16=50\", \"8+56=64\", \"20+53=73\", \"10+29=39\", \"82-8=74\", \"46+45=91\", \"42+6=48\", \"98-88=10\", \"29+29=58\", \"28-12=16\", \"9+85=94\", \"60+16=76\", \"39-2=37\", \"89-19=70\", \"1+7=8\", \"10+63=73\", \"82-48=34\", \"86-7=79\", \"88-56=32\", \"8+63=71\", \"35+21=56\", \"4+25=29\", \"43-34=9\", \"18+77=95\", \"49+14=63\", \"79-50=29\", \"15+80=95\", \"8+54=62\", \"26+42=68\", \"34-28=6\", \"3+29=32\", \"65-29=36\", \"30+28=58\", \"69-14=55\", \"13+36=49\", \"88-56=32\", \"84-28=56\", \"61-24=37\", \"50+21=71\", \"34+21=55\", \"59-18=41\", \"86-52=34\", \"80-20=60\", \"3-2=1\", \"81-35=46\", \"45+38=83\", \"5+28=33\", \"87+12=99\", \"58+33=91\", \"37-4=33\", \"41+55=96\", \"32-32=0\", \"19+64=83\", \"43-2=41\", \"50+32=82\", \"4+22=26\", \"86+4=90\", \"88-53=35\", \"71-5=66\", \"34+54=88\", \"14+65=79\", \"38-2=36\", \"79-22=57\", \"40+47=87\", \"82-4=78\", \"67-32=35\", \"31-13=18\", \"27+31=58\", \"52-45=7\", \"26+56=82\")\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif ($rowCount * $colCount -ne $oldValues.Count) {\n    throw \"Expected $($oldValues.Count) table cells, found $($rowCount * $colCount)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        # Strip the trailing end-of-cell marker(s) (cr + cell-mark) that\n        # Range.Text always carries for table cells.\n        $actual = $cellRange.Text -replace \"[\\x07\\x0d]+$\", \"\"\n        $expected = $oldValues[$i]\n        if ($actual -ne $expected) {\n            throw \"Cell $i (row $r, col $c): expected '$expected' but found '$actual'\"\n        }\n        $cellRange.Text = $newValues[$i]\n        $i++\n    }\n}\n\n$d.Save()\n"}
